$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 31.63742133333333
$ws.Range("N2").Value = 94.91226399999999
$ws.Range("O2").Value = 0.3438737102674882
$ws.Range("P2").Value = 0.3438737102674882
$ws.Range("Q2").Value = 7.373459599175111
$ws.Range("R2").Value = 66.361136392576
$ws.Range("S2").Value = 0.2446465000825865
$ws.Range("T2").Value = 0.2446465000825865
$ws.Range("O3").Value = 0.1676771690762499
$ws.Range("P3").Value = 0.1676771690762499
$ws.Range("S3").Value = 0.1192927267581789
$ws.Range("T3").Value = 0.1192927267581789
$ws.Range("M4").Value = 14.51771933333333
$ws.Range("N4").Value = 43.553158
$ws.Range("O4").Value = 0.1577961098402008
$ws.Range("P4").Value = 0.1577961098402008
$ws.Range("Q4").Value = 3.383519024785778
$ws.Range("R4").Value = 30.451671223072
$ws.Range("S4").Value = 0.1122629175956008
$ws.Range("T4").Value = 0.1122629175956008
$ws.Range("M5").Value = 30.421077
$ws.Range("N5").Value = 91.263231
$ws.Range("O5").Value = 0.3306530108160611
$ws.Range("P5").Value = 0.3306530108160611
$ws.Range("Q5").Value = 7.089976767056
$ws.Range("R5").Value = 63.80979090350401
$ws.Range("S5").Value = 0.2352407276933002
$ws.Range("T5").Value = 0.2352407276933002
$ws.Range("M6").Value = 31.63742133333333
$ws.Range("N6").Value = 94.91226399999999
$ws.Range("O6").Value = 0.3438737102674882
$ws.Range("P6").Value = 0.3438737102674882
$ws.Range("Q6").Value = 2.990632709604444
$ws.Range("R6").Value = 26.91569438643999
$ws.Range("S6").Value = 0.09922721018490169
$ws.Range("T6").Value = 0.09922721018490169
$ws.Range("O7").Value = 0.1676771690762499
$ws.Range("P7").Value = 0.1676771690762499
$ws.Range("S7").Value = 0.04838444231807101
$ws.Range("T7").Value = 0.04838444231807101
$ws.Range("M8").Value = 14.51771933333333
$ws.Range("N8").Value = 43.553158
$ws.Range("O8").Value = 0.1577961098402008
$ws.Range("P8").Value = 0.1577961098402008
$ws.Range("Q8").Value = 1.372335812381111
$ws.Range("R8").Value = 12.35102231143
$ws.Range("S8").Value = 0.04553319224460005
$ws.Range("T8").Value = 0.04553319224460005
$ws.Range("M9").Value = 30.421077
$ws.Range("N9").Value = 91.263231
$ws.Range("O9").Value = 0.3306530108160611
$ws.Range("P9").Value = 0.3306530108160611
$ws.Range("Q9").Value = 2.875653707015
$ws.Range("R9").Value = 25.880883363135
$ws.Range("S9").Value = 0.095412283122761
$ws.Range("T9").Value = 0.095412283122761
